# Update "想去人数" (want-to-go count) values in F column for sheets
# "展览" and "全部类型" to reflect newly scraped data.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F, identical target values for both sheets.
$updates = @{
    5  = 33
    8  = 1795
    9  = 6432
    10 = 489
    13 = 106
    16 = 6920
    18 = 1299
    19 = 139
    20 = 118
    22 = 109
    23 = 281
    26 = 11
    27 = 102
    29 = 396
    30 = 593
    33 = 50
    37 = 64
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
